# petty-cashBook-2021.xlsx : "Update 1-Feb-2021, midday update."
# Clears the daily transaction entries on Sheet1 (rows 3-47), leaving only
# the running-balance column E (and the pre-existing blank helper cells in
# F/G), updates the opening balance in E2, moves the date in A3 forward to
# 1-Feb-2021 (serial 44228), and refreshes the frozen-pane/selection to the
# top of the now-emptied ledger.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Opening balance (E2) -------------------------------------------------
$ws.Range("E2").Value = 433525

# --- Row 3: keep A3 (new date) and B3, drop D3 ----------------------------
$ws.Range("A3").Value = 44228
$ws.Range("D3").Clear()

# --- Rows 4-47: clear every A/B/C/D entry, leave E's running formula ------
$clearAddr = "B4,D4,B5,C5,B6,D6,B7,D7,B8,C8,B9,C9,B10,D10,A11,B11,D11,B12,C12,B13,D13,B14,D14,B15,D15,B16,D16,B17,D17,B18,C18,B19,C19,B20,D20,A21,B21,D21,B22,C22,B23,D23,B24,D24,B25,D25,B26,C26,B27,D27,B28,D28,A29,B29,D29,B30,C30,B31,D31,B32,D32,B33,D33,B34,D34,B35,C35,B36,D36,B37,C37,B38,D38,A39,B39,D39,B40,C40,B41,D41,B42,D42,B43,C43,B44,C44,B45,D45,A46,B46,D46,B47,D47"
$clearRange = $ws.Range($clearAddr)
foreach ($area in $clearRange.Areas) {
    $area.Clear()
}

# --- Frozen pane / active selection moves back to the top of the sheet ----
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("A4").Select() | Out-Null
